$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ts = "2026-02-21 06:13:14"

# Row 2
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 3.55
$ws.Range("H2").Value = 2.38
$ws.Range("I2").Value = 2.74
$ws.Range("P2").Value = 1.7
$ws.Range("Q2").Value = 2.14
$ws.Range("BH2").Value = $ts

# Row 3
$ws.Range("BH3").Value = $ts

# Row 4
$ws.Range("G4").Value = 7.4
$ws.Range("H4").Value = 1.59
$ws.Range("I4").Value = 1.67
$ws.Range("K4").Value = 4.5
$ws.Range("N4").Value = 3.6
$ws.Range("P4").Value = 1.91
$ws.Range("Q4").Value = 1.94
$ws.Range("R4").Value = 1.34
$ws.Range("S4").Value = 3.35
$ws.Range("V4").Value = 2.46
$ws.Range("W4").Value = 1.16
$ws.Range("Z4").Value = 9.4
$ws.Range("AA4").Value = 16
$ws.Range("AJ4").Value = 230
$ws.Range("AX4").Value = 7.6
$ws.Range("BA4").Value = 7
$ws.Range("BB4").Value = 8.2
$ws.Range("BC4").Value = 8
$ws.Range("BD4").Value = 8
$ws.Range("BE4").Value = 8.2
$ws.Range("BF4").Value = 8
$ws.Range("BH4").Value = $ts

# Row 5
$ws.Range("F5").Value = 3.3
$ws.Range("K5").Value = 5.6
$ws.Range("BH5").Value = $ts

# Row 6
$ws.Range("H6").Value = 5.8
$ws.Range("Z6").Value = 48
$ws.Range("BH6").Value = $ts

# Row 7
$ws.Range("BH7").Value = $ts

# Row 8
$ws.Range("BH8").Value = $ts

# Row 9
$ws.Range("I9").Value = 36
$ws.Range("P9").Value = 2.9
$ws.Range("Q9").Value = 1.43
$ws.Range("BH9").Value = $ts

# Row 10
$ws.Range("BH10").Value = $ts

# Row 11
$ws.Range("F11").Value = 1.45
$ws.Range("G11").Value = 1.63
$ws.Range("H11").Value = 2.66
$ws.Range("I11").Value = 18
$ws.Range("BH11").Value = $ts

# Row 12
$ws.Range("P12").Value = 2.2
$ws.Range("Q12").Value = 1.48
$ws.Range("BH12").Value = $ts

# Row 13
$ws.Range("F13").Value = 1.55
$ws.Range("G13").Value = 1.89
$ws.Range("H13").Value = 2.16
$ws.Range("I13").Value = 14
$ws.Range("J13").Value = 3.6
$ws.Range("P13").Value = 1.96
$ws.Range("Q13").Value = 1.01
$ws.Range("BH13").Value = $ts

# Row 14
$ws.Range("F14").Value = 1.39
$ws.Range("G14").Value = 1.64
$ws.Range("H14").Value = 2.64
$ws.Range("J14").Value = 2.56
$ws.Range("P14").Value = 2.46
$ws.Range("Q14").Value = 1.57
$ws.Range("BH14").Value = $ts

# Row 15
$ws.Range("BH15").Value = $ts

# Row 16
$ws.Range("F16").Value = 1.46
$ws.Range("I16").Value = 14.5
$ws.Range("BH16").Value = $ts

# Row 17
$ws.Range("J17").Value = 3.45
$ws.Range("K17").Value = 6.6
$ws.Range("BH17").Value = $ts

# Row 18
$ws.Range("BH18").Value = $ts

# Row 19
$ws.Range("P19").Value = 1.71
$ws.Range("Q19").Value = 2.16
$ws.Range("BH19").Value = $ts

# Row 20
$ws.Range("Q20").Value = 2.22
$ws.Range("BH20").Value = $ts

# Row 21
$ws.Range("F21").Value = 3.9
$ws.Range("H21").Value = 2
$ws.Range("I21").Value = 2.02
$ws.Range("U21").Value = 2.34
$ws.Range("Z21").Value = 13.5
$ws.Range("AE21").Value = 19
$ws.Range("BE21").Value = 60
$ws.Range("BG21").Value = 10
$ws.Range("BH21").Value = $ts

# Row 22
$ws.Range("BH22").Value = $ts

# Row 23
$ws.Range("BH23").Value = $ts

# Row 24
$ws.Range("H24").Value = 2.02
$ws.Range("BH24").Value = $ts

# Row 25
$ws.Range("K25").Value = 4.9
$ws.Range("BH25").Value = $ts

# Row 26
$ws.Range("BH26").Value = $ts
